$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,35
$row2[0,0] = 6.701317548751831
$row2[0,1] = 4.007112400399314
$row2[0,2] = 4.846241001729612
$row2[0,3] = 1.857068498406494
$row2[0,4] = 2.571770357066079
$row2[0,5] = 6.826751655340194
$row2[0,6] = 4.262485955442701
$row2[0,7] = 7.25578863116411
$row2[0,8] = 8.772927737898296
$row2[0,9] = 2.716920540332794
$row2[0,10] = 1.400512971516166
$row2[0,11] = 4.255518593390783
$row2[0,12] = 1.652874019410875
$row2[0,13] = 2.745673475787044
$row2[0,14] = 1.286560085202966
$row2[0,15] = 2.187896425170558
$row2[0,16] = 3.877484217818294
$row2[0,17] = 4.641423268960073
$row2[0,18] = 2.325519619243486
$row2[0,19] = 2.999435407774789
$row2[0,20] = 2.286064408042214
$row2[0,21] = 3.565319912774222
$row2[0,22] = 2.332739334819572
$row2[0,23] = 2.064780284251486
$row2[0,24] = 3.458527051976749
$row2[0,25] = 4.187849310587864
$row2[0,26] = 3.969698128917
$row2[0,27] = 5.486767811700702
$row2[0,28] = 3.20717050595717
$row2[0,29] = 3.590096049010754
$row2[0,30] = 2.401061882397958
$row2[0,31] = 3.522670024739845
$row2[0,32] = 4.2858019420611
$row2[0,33] = 3.37664715880933
$row2[0,34] = 2.674770631880131
$ws.Range("B2:AJ2").Value = $row2

$row3 = New-Object "object[,]" 1,35
$row3[0,0] = 6.481727587667324
$row3[0,1] = 3.055776136932819
$row3[0,2] = 6.176818554766428
$row3[0,3] = 1.351747352290463
$row3[0,4] = 1.936404876924423
$row3[0,5] = 6.639104261275435
$row3[0,6] = 3.623669873017528
$row3[0,7] = 8.545512778148753
$row3[0,8] = 8.645533031577967
$row3[0,9] = 2.007581665363414
$row3[0,10] = 1.886893240317852
$row3[0,11] = 5.264686003089454
$row3[0,12] = 0.7468134664451471
$row3[0,13] = 2.683667116319433
$row3[0,14] = 0.758544988209353
$row3[0,15] = 1.018857601710132
$row3[0,16] = 3.646705982603412
$row3[0,17] = 5.479626916174669
$row3[0,18] = 1.4522817778268
$row3[0,19] = 3.074376418877727
$row3[0,20] = 1.292851643797625
$row3[0,21] = 4.260537832544635
$row3[0,22] = 1.715616741341594
$row3[0,23] = 0.9855292152220859
$row3[0,24] = 2.89752618608739
$row3[0,25] = 5.002384557557412
$row3[0,26] = 2.732995512268329
$row3[0,27] = 5.156478886329017
$row3[0,28] = 1.43483308112219
$row3[0,29] = 2.641751186490346
$row3[0,30] = 1.352026371295283
$row3[0,31] = 2.484473375720245
$row3[0,32] = 3.264873885301804
$row3[0,33] = 2.373743902237176
$row3[0,34] = 1.731576289011944
$ws.Range("B3:AJ3").Value = $row3

$row4 = New-Object "object[,]" 1,35
$row4[0,0] = 9.363990514991952
$row4[0,1] = 5.073515629255379
$row4[0,2] = 7.940542589686808
$row4[0,3] = 2.289949301863005
$row4[0,4] = 3.209031480149534
$row4[0,5] = 13.03925278769952
$row4[0,6] = 7.292786992409183
$row4[0,7] = 11.45816851240514
$row4[0,8] = 12.42891164859454
$row4[0,9] = 3.449012534444893
$row4[0,10] = 2.349851630155028
$row4[0,11] = 7.252696326599764
$row4[0,12] = 1.819444848703198
$row4[0,13] = 3.839373910789983
$row4[0,14] = 1.493528490513319
$row4[0,15] = 2.413495717799522
$row4[0,16] = 5.322907887892244
$row4[0,17] = 7.340220524186193
$row4[0,18] = 2.74174467478176
$row4[0,19] = 4.449568394438519
$row4[0,20] = 2.632099186978779
$row4[0,21] = 5.555509760856201
$row4[0,22] = 2.895688796708999
$row4[0,23] = 2.287921645574845
$row4[0,24] = 4.511880690833604
$row4[0,25] = 6.792405244418731
$row4[0,26] = 4.833587954444695
$row4[0,27] = 7.818276106376094
$row4[0,28] = 3.518823719637809
$row4[0,29] = 4.45731297672104
$row4[0,30] = 2.755553206124374
$row4[0,31] = 4.310662600791596
$row4[0,32] = 5.387717491990499
$row4[0,33] = 4.273328007550294
$row4[0,34] = 3.212860511624751
$ws.Range("B4:AJ4").Value = $row4

$row5 = New-Object "object[,]" 1,35
$row5[0,0] = 0.1636363636363636
$row5[0,1] = 0.3148148148148148
$row5[0,2] = 0.2407407407407407
$row5[0,3] = 0.7192982456140351
$row5[0,4] = 0.4912280701754386
$row5[0,5] = 0.45
$row5[0,6] = 0.3809523809523809
$row5[0,7] = 0.3076923076923077
$row5[0,8] = 0.2222222222222222
$row5[0,9] = 0.5
$row5[0,10] = 0.8928571428571429
$row5[0,11] = 0.5333333333333333
$row5[0,12] = 0.6851851851851852
$row5[0,13] = 0.4821428571428572
$row5[0,14] = 0.7857142857142857
$row5[0,15] = 0.4821428571428572
$row5[0,16] = 0.4285714285714285
$row5[0,17] = 0.25
$row5[0,18] = 0.5535714285714286
$row5[0,19] = 0.4489795918367347
$row5[0,20] = 0.4727272727272727
$row5[0,21] = 0.3035714285714285
$row5[0,22] = 0.4464285714285715
$row5[0,23] = 0.5357142857142857
$row5[0,24] = 0.3928571428571428
$row5[0,25] = 0.2653061224489796
$row5[0,26] = 0.2
$row5[0,27] = 0.25
$row5[0,28] = 0.2181818181818182
$row5[0,29] = 0.1964285714285714
$row5[0,30] = 0.4642857142857143
$row5[0,31] = 0.1607142857142857
$row5[0,32] = 0.2857142857142857
$row5[0,33] = 0.3695652173913043
$row5[0,34] = 0.4339622641509434
$ws.Range("B5:AJ5").Value = $row5

$row6 = New-Object "object[,]" 1,35
$row6[0,0] = 0.1818181818181818
$row6[0,1] = 0.4074074074074074
$row6[0,2] = 0.2962962962962963
$row6[0,3] = 0.7894736842105263
$row6[0,4] = 0.5964912280701754
$row6[0,5] = 0.45
$row6[0,6] = 0.5238095238095238
$row6[0,7] = 0.4230769230769231
$row6[0,8] = 0.3148148148148148
$row6[0,9] = 0.6
$row6[0,10] = 0.9642857142857143
$row6[0,11] = 0.7111111111111111
$row6[0,12] = 0.8148148148148148
$row6[0,13] = 0.5714285714285714
$row6[0,14] = 0.9464285714285714
$row6[0,15] = 0.6785714285714286
$row6[0,16] = 0.5178571428571429
$row6[0,17] = 0.4038461538461539
$row6[0,18] = 0.6607142857142857
$row6[0,19] = 0.5306122448979592
$row6[0,20] = 0.6363636363636364
$row6[0,21] = 0.5
$row6[0,22] = 0.625
$row6[0,23] = 0.6964285714285714
$row6[0,24] = 0.5535714285714286
$row6[0,25] = 0.4081632653061225
$row6[0,26] = 0.2727272727272727
$row6[0,27] = 0.3541666666666667
$row6[0,28] = 0.3090909090909091
$row6[0,29] = 0.3214285714285715
$row6[0,30] = 0.625
$row6[0,31] = 0.2678571428571428
$row6[0,32] = 0.375
$row6[0,33] = 0.4565217391304348
$row6[0,34] = 0.5471698113207547
$ws.Range("B6:AJ6").Value = $row6

$row7 = New-Object "object[,]" 1,35
$row7[0,0] = 0.2727272727272727
$row7[0,1] = 0.4814814814814815
$row7[0,2] = 0.4444444444444444
$row7[0,3] = 0.8421052631578947
$row7[0,4] = 0.7017543859649122
$row7[0,5] = 0.6
$row7[0,6] = 0.5238095238095238
$row7[0,7] = 0.4807692307692308
$row7[0,8] = 0.3703703703703703
$row7[0,9] = 0.64
$row7[0,10] = 0.9821428571428571
$row7[0,11] = 0.7555555555555555
$row7[0,12] = 0.9629629629629629
$row7[0,13] = 0.6607142857142857
$row7[0,14] = 0.9821428571428571
$row7[0,15] = 0.7678571428571429
$row7[0,16] = 0.5892857142857143
$row7[0,17] = 0.5384615384615384
$row7[0,18] = 0.7678571428571429
$row7[0,19] = 0.6530612244897959
$row7[0,20] = 0.7454545454545455
$row7[0,21] = 0.6607142857142857
$row7[0,22] = 0.7321428571428571
$row7[0,23] = 0.8214285714285714
$row7[0,24] = 0.625
$row7[0,25] = 0.6122448979591837
$row7[0,26] = 0.3818181818181818
$row7[0,27] = 0.5
$row7[0,28] = 0.4363636363636363
$row7[0,29] = 0.4285714285714285
$row7[0,30] = 0.7142857142857143
$row7[0,31] = 0.4642857142857143
$row7[0,32] = 0.4821428571428572
$row7[0,33] = 0.5217391304347826
$row7[0,34] = 0.660377358490566
$ws.Range("B7:AJ7").Value = $row7

$row8 = New-Object "object[,]" 1,35
$row8[0,0] = 0.3818181818181818
$row8[0,1] = 0.6111111111111112
$row8[0,2] = 0.5740740740740741
$row8[0,3] = 0.8947368421052632
$row8[0,4] = 0.8421052631578947
$row8[0,5] = 0.8
$row8[0,6] = 0.6666666666666666
$row8[0,7] = 0.6153846153846154
$row8[0,8] = 0.5185185185185185
$row8[0,9] = 0.8
$row8[0,10] = 0.9821428571428571
$row8[0,11] = 0.8222222222222222
$row8[0,12] = 1
$row8[0,13] = 0.7857142857142857
$row8[0,14] = 1
$row8[0,15] = 0.9464285714285714
$row8[0,16] = 0.6785714285714286
$row8[0,17] = 0.6730769230769231
$row8[0,18] = 0.8571428571428571
$row8[0,19] = 0.8163265306122449
$row8[0,20] = 0.8909090909090909
$row8[0,21] = 0.8571428571428571
$row8[0,22] = 0.9285714285714286
$row8[0,23] = 0.9464285714285714
$row8[0,24] = 0.7321428571428571
$row8[0,25] = 0.7551020408163265
$row8[0,26] = 0.6727272727272727
$row8[0,27] = 0.5625
$row8[0,28] = 0.6909090909090909
$row8[0,29] = 0.6607142857142857
$row8[0,30] = 0.875
$row8[0,31] = 0.8214285714285714
$row8[0,32] = 0.6785714285714286
$row8[0,33] = 0.7173913043478261
$row8[0,34] = 0.8490566037735849
$ws.Range("B8:AJ8").Value = $row8
